$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number and week date range) ---
$ws.Range("A8").Value = "Volume 31   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/10/2024  Through  6/16/2024"

# --- Crime statistics table updates (rows 15-31) ---

# Row 15
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("M15").Copy()
$ws.Range("C15:E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 20
$ws.Range("I16").Value = 41
$ws.Range("J16").Value = 32
$ws.Range("K16").Value = 28.125
$ws.Range("L16").Value = 41.379310344827

# Row 17
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 80
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 5.882352941176
$ws.Range("I17").Value = 101
$ws.Range("J17").Value = 108
$ws.Range("K17").Value = -6.481481481481
$ws.Range("L17").Value = -12.931034482758

# Row 18
$ws.Range("D18").Value = "'0"
$ws.Range("E18").Value = "'***.*"
$ws.Range("M18").Copy()
$ws.Range("D18:E18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C18").Value = 2
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 125
$ws.Range("I18").Value = 33
$ws.Range("K18").Value = 50
$ws.Range("L18").Value = -13.157894736842

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -46.666666666666
$ws.Range("I19").Value = 162
$ws.Range("J19").Value = 179
$ws.Range("K19").Value = -9.497206703910
$ws.Range("L19").Value = -25.345622119815

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 13
$ws.Range("H20").Value = 85.714285714285
$ws.Range("I20").Value = 35
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = -30
$ws.Range("L20").Value = -2.777777777777

# Row 21
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 44.444444444444
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = -8.860759493670
$ws.Range("I21").Value = 379
$ws.Range("J21").Value = 396
$ws.Range("K21").Value = -4.292929292929
$ws.Range("L21").Value = -14.058956916099

# Row 23
$ws.Range("G23").Value = "'0"
$ws.Range("H23").Value = "'***.*"
$ws.Range("M23").Copy()
$ws.Range("G23:H23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I23").Value = 4
$ws.Range("K23").Value = -60
$ws.Range("L23").Value = -71.428571428571

# Row 24
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -19.444444444444
$ws.Range("F24").Value = 110
$ws.Range("G24").Value = 125
$ws.Range("H24").Value = -12
$ws.Range("I24").Value = 686
$ws.Range("J24").Value = 614
$ws.Range("K24").Value = 11.726384364820
$ws.Range("L24").Value = 13.764510779436

# Row 25
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = 8.695652173913
$ws.Range("F25").Value = 76
$ws.Range("G25").Value = 68
$ws.Range("H25").Value = 11.764705882352
$ws.Range("I25").Value = 476
$ws.Range("J25").Value = 369
$ws.Range("K25").Value = 28.997289972899
$ws.Range("L25").Value = 57.095709570957

# Row 26
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 31.034482758620
$ws.Range("I26").Value = 273
$ws.Range("J26").Value = 249
$ws.Range("K26").Value = 9.638554216867
$ws.Range("L26").Value = 24.090909090909

# Row 27
$ws.Range("C27").Value = "'0"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "'***.*"
$ws.Range("M27").Copy()
$ws.Range("C27:E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 300
$ws.Range("L27").Value = 87.5

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -28.571428571428
$ws.Range("I28").Value = 28
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -22.222222222222
$ws.Range("L28").Value = 7.692307692307

# Row 29
$ws.Range("L29").Value = -80

# Row 30
$ws.Range("L30").Value = -80

# Row 31
$ws.Range("C31").Value = "'0"
$ws.Range("M31").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = -50
$ws.Range("J31").Value = 3
$ws.Range("K31").Value = 66.666666666666

